$wb = $excel.ActiveWorkbook

# --- 1. Insert a new "Reference Pinouts" sheet before the first sheet ---
$firstSheet = $wb.Worksheets.Item(1)
$refSheet = $wb.Worksheets.Add($firstSheet)
$refSheet.Name = "Reference Pinouts"

# Row 1: PJRC forum thread comparing Teensy 3.2/3.5 pin compatibility
$refSheet.Hyperlinks.Add($refSheet.Range("A1"), "https://forum.pjrc.com/threads/55568-Teensy-3-2-3-5-pin-compatibility") | Out-Null
$refSheet.Range("A1").Value = "https://forum.pjrc.com/threads/55568-Teensy-3-2-3-5-pin-compatibility"
$refSheet.Range("A1").Style = "Hyperlink"

# Row 2: Google Sheets htmlview link (trailing '#' kept in the visible text)
$refSheet.Hyperlinks.Add($refSheet.Range("A2"), "https://docs.google.com/spreadsheets/u/0/d/1LSi0c17iqtvpKuNSYksMG306_FpWdJcniSRR6aGNNYQ/htmlview#", "", "", "https://docs.google.com/spreadsheets/u/0/d/1LSi0c17iqtvpKuNSYksMG306_FpWdJcniSRR6aGNNYQ/htmlview") | Out-Null
$refSheet.Range("A2").Value = "https://docs.google.com/spreadsheets/u/0/d/1LSi0c17iqtvpKuNSYksMG306_FpWdJcniSRR6aGNNYQ/htmlview#"
$refSheet.Range("A2").Style = "Hyperlink"

# Row 3: Google Sheets edit link with a gid fragment (split into Address + SubAddress)
$refSheet.Hyperlinks.Add($refSheet.Range("A3"), "https://docs.google.com/spreadsheets/d/1LSi0c17iqtvpKuNSYksMG306_FpWdJcniSRR6aGNNYQ/edit", "gid=1683806103", "", "https://docs.google.com/spreadsheets/d/1LSi0c17iqtvpKuNSYksMG306_FpWdJcniSRR6aGNNYQ/edit - gid=1683806103") | Out-Null
$refSheet.Range("A3").Value = "https://docs.google.com/spreadsheets/d/1LSi0c17iqtvpKuNSYksMG306_FpWdJcniSRR6aGNNYQ/edit#gid=1683806103"
$refSheet.Range("A3").Style = "Hyperlink"

# Remembered UI selection on this sheet (B6), then move focus elsewhere
$refSheet.Activate() | Out-Null
$refSheet.Range("B6").Select() | Out-Null

# --- 2. Rename the v0.3 sheet to reflect it now also covers Teensy LC ---
$v03 = $wb.Worksheets.Item("v0.3 Teensy 3.2")
$v03.Name = "v0.3 Teensy LC and 3.2"

# --- 3. Make the v0.3 sheet the active tab, matching the saved view state ---
$v03.Activate() | Out-Null
